$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("A2").Value = "RegisterWithoutLastName"
$ws.Range("B2").Value = "Iliya"
$ws.Range("C2").Value = "Iliev"
